$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-stagger columns D-M: each column previously had its 16 data points spread
# out down staggered rows (one new row per column per year); move them back so
# every column lines up on the same contiguous rows 3-18.
$ws.Range("D3").Value = 97.1
$ws.Range("E3").Value = 106.23
$ws.Range("F3").Value = 2.71
$ws.Range("G3").Value = 39.200000000000003
$ws.Range("I3").Value = 11.1
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 11.4
$ws.Range("L3").Value = 1.6
$ws.Range("M3").Value = 9.1
$ws.Range("D4").Value = 92
$ws.Range("E4").Value = 110.43
$ws.Range("F4").Value = 2.81
$ws.Range("G4").Value = 39.299999999999997
$ws.Range("I4").Value = 11.4
$ws.Range("J4").Value = 3.6
$ws.Range("K4").Value = 11.2
$ws.Range("L4").Value = 1.7
$ws.Range("M4").Value = 8.8000000000000007
$ws.Range("D5").Value = 96.3
$ws.Range("E5").Value = 117.2
$ws.Range("F5").Value = 2.93
$ws.Range("G5").Value = 40
$ws.Range("I5").Value = 10.5
$ws.Range("J5").Value = 3.3
$ws.Range("K5").Value = 10.8
$ws.Range("L5").Value = 1.6
$ws.Range("M5").Value = 8.6
$ws.Range("D6").Value = 96.4
$ws.Range("E6").Value = 121.6
$ws.Range("F6").Value = 3.01
$ws.Range("G6").Value = 40.4
$ws.Range("I6").Value = 10.3
$ws.Range("J6").Value = 3.2
$ws.Range("K6").Value = 10.199999999999999
$ws.Range("L6").Value = 1.6
$ws.Range("M6").Value = 7.9
$ws.Range("D7").Value = 97.4
$ws.Range("E7").Value = 127.92
$ws.Range("F7").Value = 3.12
$ws.Range("G7").Value = 41
$ws.Range("I7").Value = 9.6999999999999993
$ws.Range("J7").Value = 3.6
$ws.Range("K7").Value = 10.6
$ws.Range("L7").Value = 1.6
$ws.Range("M7").Value = 8.1999999999999993
$ws.Range("D8").Value = 97.9
$ws.Range("E8").Value = 128.21
$ws.Range("F8").Value = 3.15
$ws.Range("G8").Value = 40.700000000000003
$ws.Range("I8").Value = 11.5
$ws.Range("J8").Value = 4.8
$ws.Range("K8").Value = 10.199999999999999
$ws.Range("L8").Value = 1.8
$ws.Range("M8").Value = 7.7
$ws.Range("D9").Value = 108
$ws.Range("E9").Value = 127.98
$ws.Range("F9").Value = 3.16
$ws.Range("G9").Value = 40.5
$ws.Range("I9").Value = 9.8000000000000007
$ws.Range("J9").Value = 5.0999999999999996
$ws.Range("K9").Value = 9
$ws.Range("L9").Value = 2.2000000000000002
$ws.Range("M9").Value = 5.8
$ws.Range("D10").Value = 119.5
$ws.Range("E10").Value = 137.78
$ws.Range("F10").Value = 3.32
$ws.Range("G10").Value = 41.5
$ws.Range("I10").Value = 8.9
$ws.Range("J10").Value = 5
$ws.Range("K10").Value = 8.9
$ws.Range("L10").Value = 2.9
$ws.Range("M10").Value = 4.8
$ws.Range("D11").Value = 114.8
$ws.Range("E11").Value = 139.32
$ws.Range("F11").Value = 3.44
$ws.Range("G11").Value = 40.5
$ws.Range("I11").Value = 8.6
$ws.Range("J11").Value = 4.5999999999999996
$ws.Range("K11").Value = 8.9
$ws.Range("L11").Value = 2.4
$ws.Range("M11").Value = 5.4
$ws.Range("D12").Value = 115.4
$ws.Range("E12").Value = 144.99
$ws.Range("F12").Value = 3.58
$ws.Range("G12").Value = 40.5
$ws.Range("I12").Value = 9.1
$ws.Range("J12").Value = 4.9000000000000004
$ws.Range("K12").Value = 8.6
$ws.Range("L12").Value = 2.7
$ws.Range("M12").Value = 4.5999999999999996
$ws.Range("D13").Value = 115.7
$ws.Range("E13").Value = 155.07
$ws.Range("F13").Value = 3.81
$ws.Range("G13").Value = 40.700000000000003
$ws.Range("I13").Value = 8.3000000000000007
$ws.Range("J13").Value = 5
$ws.Range("K13").Value = 8
$ws.Range("L13").Value = 3
$ws.Range("M13").Value = 3.7
$ws.Range("D14").Value = 106.7
$ws.Range("E14").Value = 158
$ws.Range("F14").Value = 3.96
$ws.Range("G14").Value = 39.9
$ws.Range("I14").Value = 7.3
$ws.Range("J14").Value = 3.7
$ws.Range("K14").Value = 8.1
$ws.Range("L14").Value = 2.2000000000000002
$ws.Range("M14").Value = 4.5999999999999996
$ws.Range("D15").Value = 104.7
$ws.Range("E15").Value = 162.74
$ws.Range("F15").Value = 4.12
$ws.Range("G15").Value = 39.5
$ws.Range("I15").Value = 8.6
$ws.Range("J15").Value = 3.9
$ws.Range("K15").Value = 8.6
$ws.Range("L15").Value = 2
$ws.Range("M15").Value = 5.3
$ws.Range("D16").Value = 109.8
$ws.Range("E16").Value = 172.66
$ws.Range("F16").Value = 4.3600000000000003
$ws.Range("G16").Value = 39.6
$ws.Range("I16").Value = 8.1
$ws.Range("J16").Value = 3.9
$ws.Range("K16").Value = 7.8
$ws.Range("L16").Value = 2.1
$ws.Range("M16").Value = 4.2
$ws.Range("D17").Value = 114
$ws.Range("E17").Value = 178.41
$ws.Range("F17").Value = 4.6100000000000003
$ws.Range("G17").Value = 38.700000000000003
$ws.Range("I17").Value = 7.9
$ws.Range("J17").Value = 4.8
$ws.Range("K17").Value = 7.5
$ws.Range("L17").Value = 2.7
$ws.Range("M17").Value = 3.3
$ws.Range("D18").Value = 121.2
$ws.Range("E18").Value = 189.74
$ws.Range("F18").Value = 4.9800000000000004
$ws.Range("G18").Value = 38.1
$ws.Range("I18").Value = 7.4
$ws.Range("J18").Value = 4.9000000000000004
$ws.Range("K18").Value = 6.5
$ws.Range("L18").Value = 2.6
$ws.Range("M18").Value = 2.2000000000000002

# The staggered layout had spilled data all the way down to row 35; clear the
# now-unused tail so the sheet's used range shrinks back to row 18.
$ws.Range("A19:M35").Clear() | Out-Null

# Leave the selection where the author's last edit landed.
[void]$ws.Range("D3").Select()

Write-Host "done"
